# chore: update Sheets via scheduled runner
# Refreshes cached market-board price/profit figures (currentAveragePrice*,
# LevePrice*/LeveProfit* columns) for a handful of leve rows across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR tables.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H3").Value = 29990
$ws.Range("J3").Value = 29990
$ws.Range("L3").Value = 29990
$ws.Range("N3").Value = -30218

$ws.Range("H28").Value = 380
$ws.Range("I28").Value = 206.66667
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 206.66667
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = 278.33333
$ws.Range("N28").Value = -1870

$ws.Range("H40").Value = 304698.9
$ws.Range("I40").Value = 2998.2
$ws.Range("K40").Value = 2998.2
$ws.Range("M40").Value = -2823.2

$ws.Range("H53").Value = 281.92307
$ws.Range("I53").Value = 215.5
$ws.Range("J53").Value = 388.2
$ws.Range("K53").Value = 215.5
$ws.Range("L53").Value = 388.2
$ws.Range("M53").Value = 421.5
$ws.Range("N53").Value = -1662.2

$ws.Range("H102").Value = 29990
$ws.Range("J102").Value = 29990
$ws.Range("L102").Value = 29990
$ws.Range("N102").Value = -36480

$ws.Range("H135").Value = 2407.6
$ws.Range("I135").Value = 759.75
$ws.Range("K135").Value = 6837.75
$ws.Range("M135").Value = -4302.75

$ws.Range("H137").Value = 3156.8572
$ws.Range("I137").Value = 2849.6667
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 8549.000100000001
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -5999.000100000001
$ws.Range("N137").Value = -20100

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 5902
$ws.Range("I32").Value = 6024.675
$ws.Range("K32").Value = 6024.675
$ws.Range("M32").Value = -5737.675

$ws.Range("H45").Value = 2869.1428
$ws.Range("I45").Value = 2869.1428
$ws.Range("K45").Value = 2869.1428
$ws.Range("M45").Value = -2492.1428

$ws.Range("H61").Value = 3227.1428
$ws.Range("I61").Value = 3227.1428
$ws.Range("K61").Value = 3227.1428
$ws.Range("M61").Value = -3015.1428

$ws.Range("H122").Value = 2270.7273
$ws.Range("I122").Value = 2357.8
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 7073.400000000001
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -4623.400000000001
$ws.Range("N122").Value = -9100

$ws.Range("H132").Value = 4785.2856
$ws.Range("I132").Value = 4625
$ws.Range("K132").Value = 13875
$ws.Range("M132").Value = -11345

$ws.Range("H136").Value = 3227.1428
$ws.Range("I136").Value = 3227.1428
$ws.Range("K136").Value = 9681.428400000001
$ws.Range("M136").Value = -7131.428400000001

$ws = $wb.Worksheets("BSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()

$ws.Range("H35").Value = 45597.6
$ws.Range("J35").Value = 54497
$ws.Range("L35").Value = 54497
$ws.Range("N35").Value = -55117

$ws.Range("H49").Value = 40000
$ws.Range("J49").Value = 40000
$ws.Range("L49").Value = 40000
$ws.Range("N49").Value = -40478

$ws.Range("H99").Value = 4000
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H100").Value = 61000
$ws.Range("J100").Value = 61000
$ws.Range("L100").Value = 61000
$ws.Range("N100").Value = -63164

$ws.Range("H134").Value = 5166.8184
$ws.Range("I134").Value = 3259.4443
$ws.Range("J134").Value = 13750
$ws.Range("K134").Value = 9778.332900000001
$ws.Range("L134").Value = 41250
$ws.Range("M134").Value = -7243.332900000001
$ws.Range("N134").Value = -46320

$ws = $wb.Worksheets("CRP")
$ws.Range("H22").Value = 397.14285
$ws.Range("I22").Value = 397.14285
$ws.Range("K22").Value = 397.14285
$ws.Range("M22").Value = -47.14285000000001

$ws.Range("H25").Value = 11278.25
$ws.Range("I25").Value = 1500
$ws.Range("K25").Value = 1500
$ws.Range("M25").Value = -1326

$ws.Range("H31").Value = 927.4167
$ws.Range("I31").Value = 927.4167
$ws.Range("K31").Value = 927.4167
$ws.Range("M31").Value = -632.4167

$ws.Range("H34").Value = 927.4167
$ws.Range("I34").Value = 927.4167
$ws.Range("K34").Value = 927.4167
$ws.Range("M34").Value = -725.4167

$ws.Range("H99").Value = 2271
$ws.Range("I99").Value = 2399.5
$ws.Range("K99").Value = 2399.5
$ws.Range("M99").Value = -901.5

$ws.Range("H126").Value = 2271
$ws.Range("I126").Value = 2399.5
$ws.Range("K126").Value = 7198.5
$ws.Range("M126").Value = -4728.5

$ws = $wb.Worksheets("CUL")
$ws.Range("H40").Value = 500
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 500
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2138

$ws.Range("H129").Value = 1500
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets("GSM")
$ws.Range("H31").Value = 2730.3333
$ws.Range("I31").Value = 2730.3333
$ws.Range("K31").Value = 2730.3333
$ws.Range("M31").Value = -2438.3333

$ws.Range("H37").Value = 2730.3333
$ws.Range("I37").Value = 2730.3333
$ws.Range("K37").Value = 2730.3333
$ws.Range("M37").Value = -2453.3333

$ws.Range("H132").Value = 3136.0833
$ws.Range("I132").Value = 2229.375
$ws.Range("K132").Value = 6688.125
$ws.Range("M132").Value = -4158.125

$ws.Range("H133").Value = 90780
$ws.Range("J133").Value = 90780
$ws.Range("L133").Value = 90780
$ws.Range("N133").Value = -100900

$ws = $wb.Worksheets("LTW")
$ws.Range("H132").Value = 3824.182
$ws.Range("I132").Value = 1807.9333
$ws.Range("K132").Value = 5423.7999
$ws.Range("M132").Value = -2893.7999

$ws = $wb.Worksheets("WVR")
$ws.Range("H14").Value = 30538.46
$ws.Range("I14").Value = 5000
$ws.Range("J14").Value = 32666.666
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 32666.666
$ws.Range("M14").Value = -4832
$ws.Range("N14").Value = -33002.666

$ws.Range("H25").Value = 60000
$ws.Range("J25").Value = 60000
$ws.Range("L25").Value = 60000
$ws.Range("N25").Value = -60586

$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws.Range("H126").Value = 2738
$ws.Range("I126").Value = 2633.2
$ws.Range("K126").Value = 7899.599999999999
$ws.Range("M126").Value = -5429.599999999999
